$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving numeric-looking string values must be forced to Text format
# first, otherwise Excel auto-converts them to numbers (losing formatting like
# trailing zeros, e.g. "1.00" -> 1) or introduces float rounding artifacts.
$textCells = @("D5", "D6", "D9", "D14", "D19", "D20", "D21", "D22", "D24", "D25", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D41", "D43", "D44", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (prices, volume deltas, and the two swapped coin rows).
$ws.Range("D2").Value = "56.411.40"
$ws.Range("E2").Value = "  -3.67%  "
$ws.Range("D3").Value = "2.970.95"
$ws.Range("E3").Value = "  -6.08%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "494.65"
$ws.Range("E5").Value = "  -6.23%  "
$ws.Range("D6").Value = "134.69"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E8").Value = "  -4.83%  "
$ws.Range("D9").Value = "7.18"
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("E10").Value = "  -4.39%  "
$ws.Range("E11").Value = "  -8.05%  "
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").Value = "3.483.31"
$ws.Range("E13").Value = "  -5.93%  "
$ws.Range("D14").Value = "25.11"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").Value = "56.398.74"
$ws.Range("E15").Value = "  -3.63%  "
$ws.Range("D16").Value = "2.980.95"
$ws.Range("E16").Value = "  -5.46%  "
$ws.Range("E17").Value = "  -4.69%  "
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "12.28"
$ws.Range("E19").Value = "  -6.47%  "
$ws.Range("D20").Value = "7.74"
$ws.Range("E20").Value = "  -2.97%  "
$ws.Range("D21").Value = "324.86"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("E23").Value = "  -8.66%  "
$ws.Range("D24").Value = "61.45"
$ws.Range("E24").Value = "  -8.57%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -6.03%  "
$ws.Range("D27").Value = "0.0₃0888"
$ws.Range("E27").Value = "  -7.41%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").Value = "6.47"
$ws.Range("E29").Value = "  -6.02%  "
$ws.Range("D30").Value = "6.73"
$ws.Range("E30").Value = "  -3.04%  "
$ws.Range("D31").Value = "1.73"
$ws.Range("E31").Value = "  -7.44%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.16"
$ws.Range("E32").Value = "  -7.36%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "20.21"
$ws.Range("E33").Value = "  -6.04%  "
$ws.Range("D34").Value = "151.87"
$ws.Range("E34").Value = "  -4.57%  "
$ws.Range("D35").Value = "4.43"
$ws.Range("E35").Value = "  -9.23%  "
$ws.Range("E36").Value = "  -8.38%  "
$ws.Range("D37").Value = "5.57"
$ws.Range("E37").Value = "  -11.30%  "
$ws.Range("D38").Value = "0.0667"
$ws.Range("E38").Value = "  -3.28%  "
$ws.Range("D39").Value = "23.01"
$ws.Range("E39").Value = "  -5.36%  "
$ws.Range("D40").Value = "3.003.68"
$ws.Range("E40").Value = "  -5.97%  "
$ws.Range("D41").Value = "36.47"
$ws.Range("E41").Value = "  -9.95%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "0.637"
$ws.Range("E43").Value = "  -8.63%  "
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").Value = "  -9.80%  "
$ws.Range("D45").Value = "2.216.51"
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("E46").Value = "  -4.74%  "
$ws.Range("E47").Value = "  -10.04%  "
$ws.Range("D48").Value = "1.94"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("D49").Value = "0.0235"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  -7.46%  "
$ws.Range("D51").Value = "18.82"
$ws.Range("E51").Value = "  -9.80%  "
